# Reemplazando libreria js-pattern porque no resulto muy potente.
# Actualiza los valores de presion (B2:B4) y agrega una nueva fila (6)
# con el tag PM_IPA_CENTRIFUGADO_MARCHA y su formula de estabilidad.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Los tres registros de fermentacion ahora reportan 9.2 en vez de 3.65
$ws.Range("B2").Value = 9.1999999999999993
$ws.Range("B3").Value = 9.1999999999999993
$ws.Range("B4").Value = 9.1999999999999993

# Nueva fila 6: PM_IPA_CENTRIFUGADO_MARCHA
$ws.Range("A6").Value = "PM_IPA_CENTRIFUGADO_MARCHA"
$ws.Range("B6").Value = 4
$ws.Range("C6").Value = 1
$ws.Range("D6").Formula = '=IF(AND(B6>3,B6<7),"presion no estable","presion estable")'

# E6 debe quedar como el texto literal "=" (igual que E2:E5), no como formula.
# Se escribe primero como formula y luego se "endurece" a valor con
# PasteSpecial para que quede almacenado como texto plano.
$ws.Range("E6").Formula = '="="'
$ws.Range("E6").Copy()
$ws.Range("E6").PasteSpecial(-4163)

$ws.Range("F6").Value = "fill"

# La seleccion activa queda en G19
$ws.Range("G19").Select()
